# Update loading_percent values for the "case with 380 kV" result set.
# Only columns B, C, E, F, G, H, J, K (indices 2,3,5,6,7,8,10,11) change,
# for rows 2-25. Columns A, D, I, L, M, N, O are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2 = @{ 2=10.87158718709192; 3=4.667729413407372; 5=20.58565607656163; 6=43.64623357156355; 7=38.1587001304758; 8=16.61509528955492; 10=8.687948621700931; 11=10.20548480321394 }
    3 = @{ 2=10.60796486234777; 3=4.448602240158425; 5=20.29518443454679; 6=43.44934364049459; 7=38.21684091622702; 8=16.67230548386124; 10=8.714749334005193; 11=10.02821525823808 }
    4 = @{ 2=10.44529183209191; 3=4.307451264123833; 5=20.11961026510542; 6=43.34053284159417; 7=38.26560997654953; 8=16.71060764341449; 10=8.732287122575853; 11=9.919977452353894 }
    5 = @{ 2=10.37890274463551; 3=4.248309354763642; 5=20.04884658990262; 6=43.2992591827435; 7=38.28875097754684; 8=16.72701257966661; 10=8.739706250808604; 11=9.876081656116957 }
    6 = @{ 2=10.36787607050932; 3=4.238392224738659; 5=20.03714618898142; 6=43.29259178192333; 7=38.29279024899593; 8=16.72978465884336; 10=8.740954647322823; 11=9.868807342387598 }
    7 = @{ 2=10.44439674341286; 3=4.306660164876272; 5=20.11865263345566; 6=43.33996375401579; 7=38.26590886261654; 8=16.71082566364421; 10=8.732386076519946; 11=9.91938452107896 }
    8 = @{ 2=10.78091940377398; 3=4.59356699887858; 5=20.48497713050377; 6=43.57585798322037; 7=38.17602320409305; 8=16.63416131193975; 10=8.696965078246457; 11=10.14427194459969 }
    9 = @{ 2=11.43005134235696; 3=5.10236940301436; 5=21.22154739558153; 6=44.13276920507256; 7=38.10424393754484; 8=16.50909631695215; 10=8.636078907716241; 11=10.58750738121705 }
    10 = @{ 2=11.894886708939; 3=5.44188138439491; 5=21.76868147705625; 6=44.59708825301784; 7=38.11613273232017; 8=16.4327307488353; 10=8.59655743597202; 11=10.91114153546983 }
    11 = @{ 2=12.10267778595511; 3=5.588661042909658; 5=22.01781290751879; 6=44.8197096646835; 7=38.13572547310181; 8=16.40138391437333; 10=8.579706611252522; 11=11.05725421542438 }
    12 = @{ 2=12.18075655507747; 3=5.643126236632753; 5=22.11209940753659; 6=44.90559485574267; 7=38.14519296174733; 8=16.39000346371414; 10=8.573487600384007; 11=11.11236990731761 }
    13 = @{ 2=12.16396911102053; 3=5.631446067246233; 5=22.09179691110872; 6=44.88702845339677; 7=38.14306276932898; 8=16.39243262480801; 10=8.574819771248835; 11=11.10051011328169 }
    14 = @{ 2=12.10911397963378; 3=5.593164388594218; 5=22.02557149710637; 6=44.82674409234556; 7=38.13646328186372; 8=16.40043780717988; 10=8.579191722857068; 11=11.0617932079013 }
    15 = @{ 2=12.07543227878321; 3=5.569569893887757; 5=21.98499688635022; 6=44.79002265158028; 7=38.13268784144896; 8=16.40540507636071; 10=8.581890766953414; 11=11.0380485453017 }
    16 = @{ 2=11.88122632837992; 3=5.432133400018908; 5=21.75239732577466; 6=44.58276410902391; 7=38.11513850633639; 8=16.43484778003442; 10=8.597681357699036; 11=10.90156568349768 }
    17 = @{ 2=11.76108956050376; 3=5.345845640234828; 5=21.60970225690581; 6=44.45849880503717; 7=38.10801233175457; 8=16.45378027119226; 10=8.607657106892059; 11=10.81751405171427 }
    18 = @{ 2=11.69164951365132; 3=5.29549430344555; 5=21.52765531878837; 6=44.38810132478107; 7=38.1052483058191; 8=16.46498890812614; 10=8.613501044727521; 11=10.76906838463938 }
    19 = @{ 2=11.66808236245281; 3=5.278322890041653; 5=21.49988308652329; 6=44.36445251350028; 7=38.1045414029895; 8=16.46883871285398; 10=8.615497937282491; 11=10.75264979841458 }
    20 = @{ 2=11.77391417155092; 3=5.35510584625719; 5=21.62489019269349; 6=44.47161602633155; 7=38.10863271328054; 8=16.4517318269944; 10=8.606584186123809; 11=10.82647238702575 }
    21 = @{ 2=12.12524334882623; 3=5.604439064992679; 5=22.0450256679351; 6=44.84440856352606; 7=38.13834607049979; 8=16.39807317815535; 10=8.57790317821016; 11=11.07317153020259 }
    22 = @{ 2=12.35128524519802; 3=5.760876771564715; 5=22.3192536611012; 6=45.09725042509085; 7=38.16970574485753; 8=16.36586071661452; 10=8.560102867449663; 11=11.23313291067927 }
    23 = @{ 2=12.23099415172292; 3=5.677983333500607; 5=22.17295399430213; 6=44.96148148695822; 7=38.15187382393519; 8=16.38279104423325; 10=8.569516857287157; 11=11.1478920754642 }
    24 = @{ 2=11.76811731684311; 3=5.350921621419538; 5=21.61802375020699; 6=44.46568247327379; 7=38.10834808735537; 8=16.45265691915989; 10=8.607068914848515; 11=10.82242270337828 }
    25 = @{ 2=11.25618885027987; 3=4.970659857387458; 5=21.0208674438748; 6=43.97223726648498; 7=38.11237165614197; 8=16.5402123627178; 10=8.651633932780998; 11=10.46771574465559 }
}

foreach ($row in $values.Keys) {
    foreach ($col in $values[$row].Keys) {
        $ws.Cells.Item($row, $col).Value = $values[$row][$col]
    }
}
